# Auto commit at 2025-10-20  7:36:30.77
# Appends two new daily rows (98 and 99) to Sheet1, mirroring the existing
# per-station daily statistics pattern, and updates the active selection
# to follow the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 98: 四方坪站 (station 2) for 2025-10-19 (serial 45949) ----
$ws.Range("A98").Value = 45949
$ws.Range("B98").Value = "四方坪站"
$ws.Range("C98").Formula = "=16400/126"
$ws.Range("D98").Formula = "=C98/(24*60)"
$ws.Range("E98").Formula = "=8921.04/126"
$ws.Range("F98").Formula = "=3008.65/126"
$ws.Range("G98").Formula = "=8921.04/(16400/60)"
$ws.Range("H98").Formula = "=369/126"

# ---- Row 99: 高岭站 (station 3) for 2025-10-19 (serial 45949) ----
$ws.Range("A99").Value = 45949
$ws.Range("B99").Value = "高岭站"
$ws.Range("C99").Formula = "=6291/36"
$ws.Range("D99").Formula = "=C99/(24*60)"
$ws.Range("E99").Formula = "=4269.53/36"
$ws.Range("F99").Formula = "=1039.3/36"
$ws.Range("G99").Formula = "=4269.53/(6291/60)"
$ws.Range("H99").Formula = "=150/36"

# Move the active selection to follow the newly added data (mirrors the
# original workbook's convention of selecting the row below the last one).
$ws.Range("J100").Select() | Out-Null
